$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the refreshed coin data (prices/1h-volumes updated by the scraping Action;
# rows 8-17 also shifted up one coin-ranking position, so names/links/values move).
# D/E columns hold numbers-as-text (e.g. "300.05", "-1.16%"); prefix with a quote
# (Excel's "treat as text" marker) so COM keeps them as text instead of coercing
# them into numeric/percentage values.
$q = "'"

$ws.Range('D2').Value = $q + '300.05'
$ws.Range('E2').Value = $q + '-1.16%'
$ws.Range('D3').Value = $q + '36.52'
$ws.Range('E3').Value = $q + '2.50%'
$ws.Range('D4').Value = $q + '4.982'
$ws.Range('E4').Value = $q + '-2.04%'
$ws.Range('D5').Value = $q + '0.07675'
$ws.Range('E5').Value = $q + '-1.44%'
$ws.Range('D6').Value = $q + '2.051'
$ws.Range('E6').Value = $q + '-9.36%'
$ws.Range('D7').Value = $q + '7.910'
$ws.Range('E7').Value = $q + '-2.08%'
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D8').Value = $q + '0.9198'
$ws.Range('E8').Value = $q + '-0.95%'
$ws.Range('B9').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C9').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D9').Value = $q + '0.09679'
$ws.Range('E9').Value = $q + '4.70%'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').Value = $q + '0.1859'
$ws.Range('E10').Value = $q + '1.38%'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').Value = $q + '0.08489'
$ws.Range('E11').Value = $q + '-0.49%'
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').Value = $q + '0.03513'
$ws.Range('E12').Value = $q + '-6.31%'
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').Value = $q + '0.09951'
$ws.Range('E13').Value = $q + '0.14%'
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').Value = $q + '0.001482'
$ws.Range('E14').Value = $q + '0.22%'
$ws.Range('B15').Value = 'TigerCash'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D15').Value = $q + '0.005636'
$ws.Range('E15').Value = $q + '-1.67%'
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D16').Value = $q + '3.462'
$ws.Range('E16').Value = $q + '-0.69%'
$ws.Range('B17').Value = 'GateToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D17').Value = $q + '4.021'
$ws.Range('E17').Value = $q + '-0.72%'
$ws.Range('D18').Value = $q + '2.288'
$ws.Range('E18').Value = $q + '4.76%'
$ws.Range('D20').Value = $q + '0.1334'
$ws.Range('E20').Value = $q + '0.88%'
$ws.Range('D21').Value = $q + '4.792'
$ws.Range('E21').Value = $q + '3.94%'
$ws.Range('E22').Value = $q + '-1.71%'
$ws.Range('D23').Value = $q + '0.04587'
$ws.Range('E23').Value = $q + '-1.90%'
$ws.Range('D24').Value = $q + '0.005085'
$ws.Range('E24').Value = $q + '12.21%'
$ws.Range('D25').Value = $q + '0.001230'
$ws.Range('E25').Value = $q + '-0.38%'
$ws.Range('E26').Value = $q + '7.40%'
$ws.Range('E39').Value = $q + '-1.05%'
$ws.Range('D40').Value = $q + '0.04608'
$ws.Range('E40').Value = $q + '-2.59%'
$ws.Range('D41').Value = $q + '0.007423'
$ws.Range('E41').Value = $q + '-7.22%'
$ws.Range('E42').Value = $q + '-1.92%'
$ws.Range('D43').Value = $q + '0.007717'
$ws.Range('E43').Value = $q + '-2.10%'
$ws.Range('E44').Value = $q + '0.59%'
$ws.Range('D45').Value = $q + '0.01038'
$ws.Range('E45').Value = $q + '7.75%'
$ws.Range('D46').Value = $q + '0.00006282'
$ws.Range('E46').Value = $q + '1.34%'
$ws.Range('D47').Value = $q + '0.00000000750'
$ws.Range('E47').Value = $q + '-0.28%'
$ws.Range('D48').Value = $q + '0.0005798'
$ws.Range('E48').Value = $q + '-0.05%'
$ws.Range('D49').Value = $q + '35.00'
$ws.Range('E49').Value = $q + '563.02%'
$ws.Range('D50').Value = $q + '0.001999'
$ws.Range('E50').Value = $q + '-25.95%'
$ws.Range('D51').Value = $q + '0.00002101'
$ws.Range('E51').Value = $q + '-0.28%'
